$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting old rows 13-21 down to 14-22
$ws.Rows.Item(13).Insert()

# The newly inserted row 13 has no A-cell in the target; remove stray formatting
$ws.Cells.Item(13,1).Clear()

# Copy B/C formatting from row 14 (shifted old row 13) down into new row 13
$ws.Cells.Item(14,2).Copy()
$ws.Cells.Item(13,2).PasteSpecial(-4122)
$ws.Cells.Item(14,3).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10
$ws.Cells.Item(10,1).Value = 'Objetivos:'
$ws.Cells.Item(10,2).Value = 'Aprofundar os conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'
$ws.Cells.Item(10,3).Value = 'Aprofundar os conceitos técnicos fundamentais de um curso de Engenharia de Produção, tendo em vista a sua formação generalista voltada para os mais diversos tipos de sistemas de produção.'

# Row 13
$ws.Cells.Item(13,2).Value = '5840535 - Messias Borges Silva'
$ws.Cells.Item(13,3).Value = '5840535 - Messias Borges Silva'

# Row 14
$ws.Cells.Item(14,1).Value = 'Programa resumido:'
$ws.Cells.Item(14,2).Value = '1 – Planejamento e Controle da Qualidade2 – Melhoramentos da Produção3 – Desafios da produção4 – Controle da Qualidade'
$ws.Cells.Item(14,3).Value = '1 – Planejamento e Controle da Qualidade2 – Melhoramentos da Produção3 – Desafios da produção4 – Controle da Qualidade'

# Row 15
$ws.Cells.Item(15,1).Value = 'Short syllabus:'
$ws.Cells.Item(15,2).Value = '1 - Quality Planning and Control2 - Production Improvements3 - Production challenges4 - Quality Control'
$ws.Cells.Item(15,3).Value = '1 - Quality Planning and Control2 - Production Improvements3 - Production challenges4 - Quality Control'

# Row 16
$ws.Cells.Item(16,1).Value = 'Programa:'
$ws.Cells.Item(16,2).Value = '1 – Planejamento e Controle da QualidadeIntrodução. Planejamento e Controle da qualidade.2 – Melhoramentos da ProduçãoIntrodução. Medidas e melhoramentos de desempenho. Prevenção e Recuperação de falhas. Administração da Qualidade Total.3 – Desafios da produçãoIntrodução. Tipo e formas de estratégias.4 - CONTROLE DA QUALIDADEAs Sete Ferramentas da Qualidade: Diagrama de Ishikawa, Histograma, Folha de Verificação, Estratificação, Diagrama de Pareto, Diagrama de Dispersão, Gráficos de Controle. Círculos de Controle da Qualidade'
$ws.Cells.Item(16,3).Value = '1 – Planejamento e Controle da QualidadeIntrodução. Planejamento e Controle da qualidade.2 – Melhoramentos da ProduçãoIntrodução. Medidas e melhoramentos de desempenho. Prevenção e Recuperação de falhas. Administração da Qualidade Total.3 – Desafios da produçãoIntrodução. Tipo e formas de estratégias.4 - CONTROLE DA QUALIDADEAs Sete Ferramentas da Qualidade: Diagrama de Ishikawa, Histograma, Folha de Verificação, Estratificação, Diagrama de Pareto, Diagrama de Dispersão, Gráficos de Controle. Círculos de Controle da Qualidade'

# Row 17
$ws.Cells.Item(17,1).Value = 'Syllabus:'
$ws.Cells.Item(17,2).Value = '1 - Production Planning and ControlIntroduction. Planning Quality Control.2 - Production ImprovementsIntroduction. Measures and performance improvements. Prevention and Recovery of failures. Total Quality Management.3 - Production challengesIntroduction. Types and forms of strategies.4 - QUALITY CONTROLThe Seven Quality Tools: Ishikawa Diagram, Histogram, Check Sheet, Stratification, Pareto Diagram, Dispersion Diagram, Control Charts. Quality Control Circles'
$ws.Cells.Item(17,3).Value = '1 - Production Planning and ControlIntroduction. Planning Quality Control.2 - Production ImprovementsIntroduction. Measures and performance improvements. Prevention and Recovery of failures. Total Quality Management.3 - Production challengesIntroduction. Types and forms of strategies.4 - QUALITY CONTROLThe Seven Quality Tools: Ishikawa Diagram, Histogram, Check Sheet, Stratification, Pareto Diagram, Dispersion Diagram, Control Charts. Quality Control Circles'

# Row 18
$ws.Cells.Item(18,1).Value = 'Avaliação:'

# Row 19
$ws.Cells.Item(19,1).Value = 'Método:'
$ws.Cells.Item(19,2).Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Cells.Item(19,3).Value = 'Aulas Expositivas; trabalhos e seminários.'

# Row 20
$ws.Cells.Item(20,1).Value = 'Critério:'
$ws.Cells.Item(20,2).Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'
$ws.Cells.Item(20,3).Value = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'

# Row 21
$ws.Cells.Item(21,1).Value = 'Norma de recuperação:'
$ws.Cells.Item(21,2).Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'
$ws.Cells.Item(21,3).Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.'

# Row 22
$ws.Cells.Item(22,1).Value = 'Bibliografia:'
$ws.Cells.Item(22,2).Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. VENANZI, D; SILVA, O.R., Gerenciamento da Produçao e Operaçoes, LTC, 2014Textos complementares serão usados durante o curso.'
$ws.Cells.Item(22,3).Value = 'SLACK, N. et al. Administração da produção. São Paulo: Atlas, 2002. VENANZI, D; SILVA, O.R., Gerenciamento da Produçao e Operaçoes, LTC, 2014Textos complementares serão usados durante o curso.'
